# Rename speakers in column D ("Speaker") on the active sheet:
#   "HILLARY LEWIS-WOLFSEN" -> "T"
#   "STUDENT"               -> "STUDENT 1"
# Leaves all other speaker values (e.g. CYNTHIA, SAURABH, MS. Liu, JAKE, ASHANK) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D
    $val = $cell.Value2
    if ($val -eq "HILLARY LEWIS-WOLFSEN") {
        $cell.Value = "T"
    } elseif ($val -eq "STUDENT") {
        $cell.Value = "STUDENT 1"
    }
}
